$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text interpretation so the date-like strings aren't auto-converted
# into date serial numbers by Excel's input parser.
$ws.Range("C2:C4").NumberFormat = "@"

$ws.Range("C2").Value = "2024-01-15"
$ws.Range("C3").Value = "2024-02-15"
$ws.Range("C4").Value = "2024-03-15"

# Restore the original (default/general) cell formatting so only the
# shared-string values change, matching the source workbook's styling.
$ws.Range("C2:C4").ClearFormats()
